$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price column whose new value parses as a plain number need to be
# pinned to Text format first, otherwise Excel would silently reinterpret e.g.
# "1.00" as the number 1 and drop the formatting the source data relies on.
$textPriceCells = @("D5", "D6", "D7", "D10", "D14", "D18", "D19", "D20", "D22", "D23", "D24", "D30", "D35", "D37", "D40", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.162.41'
$ws.Range("E2").Value = '  -3.80%  '
$ws.Range("D3").Value = '3.323.95'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '573.86'
$ws.Range("E5").Value = '  -2.46%  '
$ws.Range("D6").Value = '181.01'
$ws.Range("E6").Value = '  -4.30%  '
$ws.Range("D7").Value = '0.615'
$ws.Range("E7").Value = '  +2.06%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -3.21%  '
$ws.Range("D10").Value = '6.65'
$ws.Range("E10").Value = '  -1.22%  '
$ws.Range("E11").Value = '  -2.73%  '
$ws.Range("D12").Value = '3.903.87'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("E13").Value = '  -1.22%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '26.66'
$ws.Range("E14").Value = '  -4.90%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '66.288.81'
$ws.Range("E15").Value = '  -3.71%  '
$ws.Range("E16").Value = '  -1.84%  '
$ws.Range("D17").Value = '3.315.84'
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").Value = '435.86'
$ws.Range("E18").Value = '  -2.62%  '
$ws.Range("D19").Value = '5.66'
$ws.Range("E19").Value = '  -2.17%  '
$ws.Range("D20").Value = '13.54'
$ws.Range("E20").Value = '  -1.60%  '
$ws.Range("E21").Value = '  -2.98%  '
$ws.Range("D22").Value = '73.18'
$ws.Range("E22").Value = '  -3.32%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '0.519'
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("E25").Value = '  -3.44%  '
$ws.Range("E26").Value = '  +1.36%  '
$ws.Range("E27").Value = '  -3.00%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("E29").Value = '  -2.89%  '
$ws.Range("D30").Value = '22.72'
$ws.Range("E30").Value = '  -1.93%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("E32").Value = '  -4.12%  '
$ws.Range("E33").Value = '  -2.46%  '
$ws.Range("E34").Value = '  -4.26%  '
$ws.Range("D35").Value = '160.69'
$ws.Range("E35").Value = '  -1.71%  '
$ws.Range("E36").Value = '  -4.56%  '
$ws.Range("D37").Value = '27.79'
$ws.Range("E38").Value = '  -6.91%  '
$ws.Range("D39").Value = '2.835.83'
$ws.Range("E39").Value = '  +5.30%  '
$ws.Range("D40").Value = '0.792'
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("E41").Value = '  -3.48%  '
$ws.Range("E42").Value = '  -4.02%  '
$ws.Range("E43").Value = '  -2.25%  '
$ws.Range("D44").Value = '0.0665'
$ws.Range("E44").Value = '  -2.54%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '2.34'
$ws.Range("E45").Value = '  -5.38%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '24.06'
$ws.Range("E46").Value = '  -3.58%  '
$ws.Range("D47").Value = '324.15'
$ws.Range("E47").Value = '  -1.95%  '
$ws.Range("D48").Value = '0.0271'
$ws.Range("E48").Value = '  -3.60%  '
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").Value = '0.972'
$ws.Range("E50").Value = '  -2.75%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '6.14'
$ws.Range("E51").Value = '  -2.38%  '
